$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(1, 1).Value = 'Datos actualizados a 1 de Abril de 2020 a las 07:20'
$ws.Cells.Item(38, 2).Value = 1771
$ws.Cells.Item(38, 3).Value = 120
$ws.Cells.Item(38, 5).Value = 1417
$ws.Cells.Item(38, 7).Value = 2
$ws.Cells.Item(38, 8).Value = 12
$ws.Cells.Item(39, 1).Value = 'India'
$ws.Cells.Item(39, 2).Value = 1590
$ws.Cells.Item(39, 3).Value = 193
$ws.Cells.Item(39, 4).Value = 148
$ws.Cells.Item(39, 5).Value = 1397
$ws.Cells.Item(39, 6).Value = 0
$ws.Cells.Item(39, 7).Value = 10
$ws.Cells.Item(39, 8).Value = 45
$ws.Cells.Item(40, 1).Value = 'Arabia Saudita'
$ws.Cells.Item(40, 2).Value = 1563
$ws.Cells.Item(40, 4).Value = 165
$ws.Cells.Item(40, 5).Value = 1388
$ws.Cells.Item(40, 6).Value = 31
$ws.Cells.Item(40, 8).Value = 10
$ws.Cells.Item(41, 1).Value = 'Indonesia'
$ws.Cells.Item(41, 2).Value = 1528
$ws.Cells.Item(41, 4).Value = 81
$ws.Cells.Item(41, 5).Value = 1311
$ws.Cells.Item(41, 6).Value = 0
$ws.Cells.Item(41, 8).Value = 136
$ws.Cells.Item(42, 1).Value = 'Finlandia'
$ws.Cells.Item(42, 2).Value = 1418
$ws.Cells.Item(42, 4).Value = 10
$ws.Cells.Item(42, 5).Value = 1391
$ws.Cells.Item(42, 6).Value = 56
$ws.Cells.Item(42, 8).Value = 17
$ws.Cells.Item(71, 2).Value = 463
$ws.Cells.Item(71, 5).Value = 414
$ws.Cells.Item(142, 1).Value = 'Guam'
$ws.Cells.Item(142, 6).Value = 0
$ws.Cells.Item(143, 1).Value = 'El Salvador'
$ws.Cells.Item(143, 6).Value = 5
$ws.Cells.Item(154, 1).Value = 'Haiti'
$ws.Cells.Item(154, 3).Value = 1
$ws.Cells.Item(154, 4).Value = 1
$ws.Cells.Item(154, 8).Value = 0
$ws.Cells.Item(155, 1).Value = 'Gabon'
$ws.Cells.Item(155, 2).Value = 16
$ws.Cells.Item(155, 8).Value = 1
$ws.Cells.Item(156, 1).Value = 'Eritrea'
$ws.Cells.Item(156, 5).Value = 15
$ws.Cells.Item(156, 8).Value = 0
$ws.Cells.Item(158, 1).Value = 'Birmania'
$ws.Cells.Item(158, 4).Value = 0
$ws.Cells.Item(158, 8).Value = 1
$ws.Cells.Item(164, 1).Value = 'Mongolia'
$ws.Cells.Item(164, 4).Value = 2
$ws.Cells.Item(164, 8).Value = 0
$ws.Cells.Item(165, 1).Value = 'Guyana'
$ws.Cells.Item(165, 4).Value = 0
$ws.Cells.Item(165, 8).Value = 2
$ws.Cells.Item(171, 1).Value = 'Siria'
$ws.Cells.Item(171, 4).Value = 0
$ws.Cells.Item(171, 8).Value = 2
$ws.Cells.Item(172, 1).Value = 'Groenlandia'
$ws.Cells.Item(172, 4).Value = 2
$ws.Cells.Item(172, 8).Value = 0
$ws.Cells.Item(173, 1).Value = 'Granada'
$ws.Cells.Item(175, 1).Value = 'Suazilandia'
$ws.Cells.Item(177, 1).Value = 'San Cristobal y Nieves'
$ws.Cells.Item(178, 1).Value = 'Mozambique'
$ws.Cells.Item(179, 1).Value = 'Guinea-Bisau'
$ws.Cells.Item(183, 1).Value = 'Angola'
$ws.Cells.Item(184, 1).Value = 'Sudan'
$ws.Cells.Item(185, 1).Value = 'San Martin (Parte Holandesa)'
$ws.Cells.Item(186, 1).Value = 'Santa Sede'
$ws.Cells.Item(190, 1).Value = 'Montserrat'
$ws.Cells.Item(191, 1).Value = 'Fiyi'
$ws.Cells.Item(192, 1).Value = 'Islas Turcas y Caicos'
$ws.Cells.Item(194, 1).Value = 'Nicaragua'
$ws.Cells.Item(194, 4).Value = 0
$ws.Cells.Item(194, 8).Value = 1
$ws.Cells.Item(195, 1).Value = 'Nepal'
$ws.Cells.Item(195, 4).Value = 1
$ws.Cells.Item(195, 8).Value = 0
$ws.Cells.Item(197, 1).Value = 'Gambia'
$ws.Cells.Item(198, 1).Value = 'Botsuana'
$ws.Cells.Item(200, 1).Value = 'Republica de Africa Central'
$ws.Cells.Item(201, 1).Value = 'Liberia'
$ws.Cells.Item(203, 1).Value = 'Burundi'
$ws.Cells.Item(203, 3).Value = 0
$ws.Cells.Item(205, 1).Value = 'Bonaire, San Eustaquio y Saba'
$ws.Cells.Item(205, 3).Value = 2
